# Generate Report for handback
# Update the "Correspond Handoff Datetime" (col D) and
# "Correspond Handback DateTime" (col G) for the 629e121d... file row
# (row 3) on both the zh-cn and de-de language report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-08 15:25:50"
$wsZhCn.Range("G3").Value = "2016-01-08 15:26:54"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-08 15:26:05"
$wsDeDe.Range("G3").Value = "2016-01-08 15:27:17"
